$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Keysborough / Sikh Temple / 1/01/21 3:00pm-5:00pm / Case visited venue -> Moorabbin data
$ws.Range("A2").Value = "Moorabbin"
$ws.Range("B2").Value = "Grape and Grain Liquor Cellars, 14/16 Station St"
$ws.Range("C2").Value = "24/12/20 1:00pm-10:00pm  28/12/20 8.05pm-8.47pm  29/12/20 12:00pm-4:00pm"
$ws.Range("D2").Value = "Case's workplace"

# Row 3: Keysborough / Sikh Temple / 1/01/21 3:00pm-6:00pm / Case visited venue -> Moorabbin data
$ws.Range("A3").Value = "Moorabbin"
$ws.Range("B3").Value = "Grape and Grain Liquor Cellars, 14/16 Station St"
$ws.Range("C3").Value = "28/12/20 8.05pm-8.47pm  29/12/20 12:00pm-4:00pm"
$ws.Range("D3").Value = "Case's workplace"

# Recalculate best-fit column widths based on the new content
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

$wb.Save()
